$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Sexo" category labels used throughout column G
# (shared strings are renamed in place so every cell referencing them updates).
$ws.Cells.Replace("Masculino", "Hombre")
$ws.Cells.Replace("Femenino", "Mujer")

# Update the active selection on the sheet to G2 (single cell).
$ws.Range("G2").Select()
